$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.08805033333333334
$ws.Range("H2").Value = 0.264151
$ws.Range("I2").Value = 0.005589762818257384
$ws.Range("J2").Value = 0.005589762818257385
$ws.Range("M2").Value = 6.045145666666667
$ws.Range("N2").Value = 18.135437
$ws.Range("O2").Value = 0.8160840232643366
$ws.Range("P2").Value = 0.8160840232643367
$ws.Range("Q2").Value = 0.5322770909985557
$ws.Range("R2").Value = 4.790493818987001
$ws.Range("S2").Value = 0.004561716129816883
$ws.Range("T2").Value = 0.004561716129816884
$ws.Range("G3").Value = 0.08805033333333334
$ws.Range("H3").Value = 0.264151
$ws.Range("I3").Value = 0.005589762818257384
$ws.Range("J3").Value = 0.005589762818257385
$ws.Range("O3").Value = 0.09212864864242169
$ws.Range("P3").Value = 0.09212864864242169
$ws.Range("Q3").Value = 0.06008936298111112
$ws.Range("R3").Value = 0.54080426683
$ws.Range("S3").Value = 0.0005149772946777075
$ws.Range("T3").Value = 0.0005149772946777075
$ws.Range("G4").Value = 0.08805033333333334
$ws.Range("H4").Value = 0.264151
$ws.Range("I4").Value = 0.005589762818257384
$ws.Range("J4").Value = 0.005589762818257385
$ws.Range("M4").Value = 0.6799149999999999
$ws.Range("N4").Value = 2.039745
$ws.Range("O4").Value = 0.09178732809324164
$ws.Range("P4").Value = 0.09178732809324165
$ws.Range("Q4").Value = 0.05986674238833333
$ws.Range("R4").Value = 0.538800681495
$ws.Range("S4").Value = 0.0005130693937627935
$ws.Range("T4").Value = 0.0005130693937627938
$ws.Range("I5").Value = 0.9470512964761942
$ws.Range("J5").Value = 0.9470512964761943
$ws.Range("M5").Value = 6.045145666666667
$ws.Range("N5").Value = 18.135437
$ws.Range("O5").Value = 0.8160840232643366
$ws.Range("P5").Value = 0.8160840232643367
$ws.Range("Q5").Value = 90.18159186795536
$ws.Range("R5").Value = 811.6343268115982
$ws.Range("S5").Value = 0.7728734322659986
$ws.Range("T5").Value = 0.7728734322659988
$ws.Range("I6").Value = 0.9470512964761942
$ws.Range("J6").Value = 0.9470512964761943
$ws.Range("O6").Value = 0.09212864864242169
$ws.Range("P6").Value = 0.09212864864242169
$ws.Range("R6").Value = 91.62631737582001
$ws.Range("S6").Value = 0.08725055613940523
$ws.Range("T6").Value = 0.08725055613940524
$ws.Range("I7").Value = 0.9470512964761942
$ws.Range("J7").Value = 0.9470512964761943
$ws.Range("M7").Value = 0.6799149999999999
$ws.Range("N7").Value = 2.039745
$ws.Range("O7").Value = 0.09178732809324164
$ws.Range("P7").Value = 0.09178732809324165
$ws.Range("Q7").Value = 10.14298420847
$ws.Range("R7").Value = 91.28685787623
$ws.Range("S7").Value = 0.0869273080707903
$ws.Range("T7").Value = 0.08692730807079033
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.7460013333333334
$ws.Range("H8").Value = 2.238004
$ws.Range("I8").Value = 0.04735894070554834
$ws.Range("J8").Value = 0.04735894070554835
$ws.Range("M8").Value = 6.045145666666667
$ws.Range("N8").Value = 18.135437
$ws.Range("O8").Value = 0.8160840232643366
$ws.Range("P8").Value = 0.8160840232643367
$ws.Range("Q8").Value = 4.509686727527557
$ws.Range("R8").Value = 40.58718054774801
$ws.Range("S8").Value = 0.03864887486852105
$ws.Range("T8").Value = 0.03864887486852106
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.7460013333333334
$ws.Range("H9").Value = 2.238004
$ws.Range("I9").Value = 0.04735894070554834
$ws.Range("J9").Value = 0.04735894070554835
$ws.Range("O9").Value = 0.09212864864242169
$ws.Range("P9").Value = 0.09212864864242169
$ws.Range("Q9").Value = 0.5091036365911112
$ws.Range("R9").Value = 4.58193272932
$ws.Range("S9").Value = 0.004363115208338746
$ws.Range("T9").Value = 0.004363115208338746
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.7460013333333334
$ws.Range("H10").Value = 2.238004
$ws.Range("I10").Value = 0.04735894070554834
$ws.Range("J10").Value = 0.04735894070554835
$ws.Range("M10").Value = 0.6799149999999999
$ws.Range("N10").Value = 2.039745
$ws.Range("O10").Value = 0.09178732809324164
$ws.Range("P10").Value = 0.09178732809324165
$ws.Range("Q10").Value = 0.5072174965533334
$ws.Range("R10").Value = 4.56495746898
$ws.Range("S10").Value = 0.004346950628688543
$ws.Range("T10").Value = 0.004346950628688544
